$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 48 first (shifts existing rows 48-67 down to 49-68); introduces no new strings.
$ws.Rows(48).Insert()

# 1. New rows 69-82 (final numbering) - introduces WoodDeckSF..SaleCondition = indices 86-99
$ws.Range("B69").Value = "WoodDeckSF"
$ws.Range("D69").Value = "Numerical"
$ws.Range("E69").Value = "Building"
$ws.Range("F69").Value = "Low"

$ws.Range("B70").Value = "OpenPorchSF"
$ws.Range("D70").Value = "Numerical"
$ws.Range("E70").Value = "Building"
$ws.Range("F70").Value = "Low"

$ws.Range("B71").Value = "EnclosedPorch"
$ws.Range("D71").Value = "Numerical"
$ws.Range("E71").Value = "Building"
$ws.Range("F71").Value = "Low"

$ws.Range("B72").Value = "3SsnPorch"
$ws.Range("D72").Value = "Numerical"
$ws.Range("E72").Value = "Building"
$ws.Range("F72").Value = "Low"

$ws.Range("B73").Value = "ScreenPorch"
$ws.Range("D73").Value = "Numerical"
$ws.Range("E73").Value = "Building"
$ws.Range("F73").Value = "Low"

$ws.Range("B74").Value = "PoolArea"
$ws.Range("D74").Value = "Numerical"
$ws.Range("E74").Value = "Space"
$ws.Range("F74").Value = "High"
$ws.Range("G74").Value = "Low"

$ws.Range("B75").Value = "PoolQC"
$ws.Range("D75").Value = "Categorical"
$ws.Range("E75").Value = "Space"
$ws.Range("F75").Value = "Low"

$ws.Range("B76").Value = "Fence"
$ws.Range("D76").Value = "Categorical"
$ws.Range("E76").Value = "Space"
$ws.Range("F76").Value = "Low"

$ws.Range("B77").Value = "MiscFeature"
$ws.Range("D77").Value = "Categorical"
$ws.Range("E77").Value = "Space"
$ws.Range("F77").Value = "Low"

$ws.Range("B78").Value = "MiscVal"
$ws.Range("D78").Value = "Numerical"
$ws.Range("E78").Value = "Space"
$ws.Range("F78").Value = "Low"

$ws.Range("B79").Value = "NoSold"
$ws.Range("D79").Value = "Numerical"
$ws.Range("E79").Value = "Space"
$ws.Range("F79").Value = "Low"

$ws.Range("B80").Value = "YrSold"
$ws.Range("D80").Value = "Numerical"
$ws.Range("E80").Value = "Space"
$ws.Range("F80").Value = "Low"

$ws.Range("B81").Value = "SaleType"
$ws.Range("D81").Value = "Categorical"
$ws.Range("E81").Value = "Space"
$ws.Range("F81").Value = "Low"

$ws.Range("B82").Value = "SaleCondition"
$ws.Range("D82").Value = "Categorical"
$ws.Range("E82").Value = "Space"
$ws.Range("F82").Value = "Low"

# 2. G-column fills (row order) - introduces Little(row20), medium(row51), HIgh(row56)
$ws.Range("G7").Value = "Low"
$ws.Range("G14").Value = "Medium"
$ws.Range("G19").Value = "High"
$ws.Range("G20").Value = "Little"
$ws.Range("G21").Value = "Low"
$ws.Range("G40").Value = "High"
$ws.Range("G51").Value = "medium"
$ws.Range("G53").Value = "Low"
$ws.Range("G56").Value = "HIgh"
$ws.Range("G65").Value = "Medium"

# 3. Row 48 content last - introduces GrLivArea
$ws.Range("B48").Value = "GrLivArea"
$ws.Range("D48").Value = "Numerical"
$ws.Range("E48").Value = "Building"
$ws.Range("F48").Value = "HIgh"
$ws.Range("G48").Value = "High"

$ws.Range("G49").Select()
